$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")

# Version: 2.1.0 -> 2.2.0-ballot
$meta.Range("B3").Value = "2.2.0-ballot"

# Date: 2025-12-18T17:25:31+00:00 -> 2025-12-19T08:32:44+00:00
$meta.Range("B8").Value = "2025-12-19T08:32:44+00:00"

# Base Definition: append version to canonical URL
$meta.Range("B18").Value = "https://hl7.fr/ig/fhir/core/StructureDefinition/fr-core-human-name|2.1.0"

# --- Elements sheet ---
$elems = $wb.Worksheets.Item("Elements")

# HumanName.extension:assemblyOrder Type(s) -> append extension version
$elems.Range("K5").Value = "Extension {humanname-assembly-order|5.2.0}`n"

# HumanName.use Binding Value Set -> append value set version
$elems.Range("Z6").Value = "http://hl7.org/fhir/ValueSet/name-use|4.0.1"

# HumanName.prefix Binding Value Set -> append value set version
$elems.Range("Z10").Value = "https://mos.esante.gouv.fr/NOS/JDV_J245-Civilite-CISIS/FHIR/JDV-J245-Civilite-CISIS|20230331120000"

# HumanName.suffix Binding Value Set -> append value set version
$elems.Range("Z11").Value = "https://mos.esante.gouv.fr/NOS/JDV_J79-CiviliteExercice-RASS/FHIR/JDV-J79-CiviliteExercice-RASS|20200424120000"

# Refresh the bestFit widths for columns K and Z to reflect the new (longer) text
$elems.Columns.Item(11).ColumnWidth = 37.5
$elems.Columns.Item(26).ColumnWidth = 91
